# worklistInfo.xlsx / worklistCalibratorsControls sheet:
# "Added around 4-6 Drop 2 defects test cases"
#
# The sheet holds three stacked test-case blocks (header band + field-header
# row + data row), separated by a blank spacer row. This change adds two
# more blank spacer rows so there is breathing room to add new test cases
# later:
#   - one new blank row before the "last1" (Assert404) block  (old row 5 -> 6)
#   - one new blank row after the blank spacer that follows the
#     "last1000001" block (old row 9 -> 11)
#
# Inserting whole rows shifts everything below down and Excel carries the
# formatting of the row immediately above into the newly inserted row,
# which reproduces the blank s4/s4/s4/s4/s2 styling used elsewhere in the
# sheet without having to touch any cell values.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("worklistCalibratorsControls")

# Insert a new blank row right before the "last1" / Assert404 block.
$ws.Rows(5).Insert()

# Insert a second new blank row right after the blank spacer that used to
# sit at row 8 (now row 9 after the first insert), pushing the Assert401
# block further down.
$ws.Rows(10).Insert()

# Reflect the refreshed view position: selection back at the top of the
# sheet.
$ws.Activate() | Out-Null
$excel.ActiveWindow.ScrollRow = 2
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("A2").Select() | Out-Null
